# Seventh Commit with Documentation
#
# - Corrects the typo'd test e-mail address (testjaga6717@gmail.com ->
#   testjaga8717@gmail.com) everywhere it appears: SignIn!C2, SignIn!C3,
#   CreateAccount!F2, CreateAccount!F3.
# - Moves the active window/tab from the "SignIn" sheet to the
#   "CreateAccount" sheet and updates each sheet's remembered selection.

$wb = $excel.ActiveWorkbook

$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsCreateAccount = $wb.Worksheets.Item("CreateAccount")

# Fix the mistyped e-mail address on both sheets.
$wsSignIn.Range("C2").Value = "testjaga8717@gmail.com"
$wsSignIn.Range("C3").Value = "testjaga8717@gmail.com"

$wsCreateAccount.Range("F2").Value = "testjaga8717@gmail.com"
$wsCreateAccount.Range("F3").Value = "testjaga8717@gmail.com"

# SignIn is no longer the active tab; its remembered selection moves to C6.
$wsSignIn.Activate()
[void]$wsSignIn.Range("C6").Select()

# CreateAccount becomes the active tab; its remembered selection moves to M10.
$wsCreateAccount.Activate()
[void]$wsCreateAccount.Range("M10").Select()
